$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8102935552597046
$ws.Range("B1").Value = 1.554062724113464
$ws.Range("C1").Value = 4.043524742126465
$ws.Range("D1").Value = 1.287571907043457
$ws.Range("E1").Value = 0.7911099791526794
